# Apply the changes described in the diff:
# 1. On "tracks_description" sheet, change the range_detection (column L) value
#    for the Density, Neutron and PhotoelectricFactor curves (rows 6, 8, 9) from
#    "manual" to "auto" -- this fixes figure display driven by matplotlib when
#    axis dimensions change.
# 2. Update the selected cell on "tracks_description" sheet to G20.
# 3. Update the selected cell on "aliases" sheet to A14.

$wb = $excel.ActiveWorkbook

$wsTracks = $wb.Worksheets.Item("tracks_description")
$wsAliases = $wb.Worksheets.Item("aliases")

# Change range_detection from "manual" to "auto" for rows 6, 8 and 9.
$wsTracks.Range("L6").Value = "auto"
$wsTracks.Range("L8").Value = "auto"
$wsTracks.Range("L9").Value = "auto"

# Update the selections to match the saved view state.
$wsTracks.Activate()
$wsTracks.Range("G20").Select()

$wsAliases.Activate()
$wsAliases.Range("A14").Select()

# Restore the original active sheet (tracks_description was tabSelected).
$wsTracks.Activate()
